# "preparing dynamic data for NL case"
# Shifts the simulation window two years later (2025 -> 2028) and swaps in a
# fresh, shorter set of dynamic plant records for conventionals, renewables,
# storages and biogas.

$wb = $excel.ActiveWorkbook

# --- times: StartTime / StopTime pushed out by 2 years (1095 days) ---
$wsTimes = $wb.Worksheets.Item("times")
$wsTimes.Range("B2").Value = 46752.99861111111
$wsTimes.Range("B3").Value = 47117.99861111111

# --- scenario_data_emlab: base year 2025 -> 2028 ---
$wsScenario = $wb.Worksheets.Item("scenario_data_emlab")
$wsScenario.Range("B1").Value = 2028

# --- conventionals: replace rows 2-15 with new data, drop rows 16-21 ---
$wsConv = $wb.Worksheets.Item("conventionals")

$wsConv.Range("A16:G21").Delete()

$convData = @(
  @(0, 20030300025, "NATURAL_GAS", 4.2, 0.61, 5667.109666666667, 5667.109666666667),
  @(1, 19980300022, "NATURAL_GAS", 4.2, 0.61, 5667.109666666667, 5667.109666666667),
  @(2, 19970300031, "NATURAL_GAS", 4.2, 0.61, 2481.7675, 2481.7675),
  @(3, 19940300026, "NATURAL_GAS", 4.2, 0.61, 7561.923333333333, 7561.923333333333),
  @(4, 19920300041, "NATURAL_GAS", 4.2, 0.61, 2481.7675, 2481.7675),
  @(5, 19900300030, "NATURAL_GAS", 4.2, 0.61, 2704.666666666667, 2704.666666666667),
  @(6, 19890300027, "NATURAL_GAS", 4.2, 0.61, 1217.633333333333, 1217.633333333333),
  @(7, 19890300033, "NATURAL_GAS", 4.2, 0.61, 7561.923333333333, 7561.923333333333),
  @(8, 19870300029, "NATURAL_GAS", 4.2, 0.61, 6926.339999999999, 6926.339999999999),
  @(9, 19850300040, "NATURAL_GAS", 4.2, 0.61, 2704.666666666667, 2704.666666666667),
  @(10, 19840300037, "NATURAL_GAS", 4.2, 0.61, 1217.633333333333, 1217.633333333333),
  @(11, 19820300039, "NATURAL_GAS", 4.2, 0.61, 6926.339999999999, 6926.339999999999),
  @(12, 19690300028, "NATURAL_GAS", 4.2, 0.61, 2940.25, 2940.25),
  @(13, 19640300038, "NATURAL_GAS", 4.2, 0.61, 2940.25, 2940.25)
)

for ($i = 0; $i -lt $convData.Count; $i++) {
    $r = $i + 2
    $row = $convData[$i]
    $wsConv.Cells.Item($r, 1).Value = $row[0]
    $wsConv.Cells.Item($r, 2).Value = $row[1]
    $wsConv.Cells.Item($r, 3).Value = $row[2]
    $wsConv.Cells.Item($r, 4).Value = $row[3]
    $wsConv.Cells.Item($r, 5).Value = $row[4]
    $wsConv.Cells.Item($r, 6).Value = $row[5]
    $wsConv.Cells.Item($r, 7).Value = $row[6]
}

# --- renewables: replace rows 2-10 with new data, drop rows 11-16 ---
$wsRenew = $wb.Worksheets.Item("renewables")

$wsRenew.Range("A11:I16").Delete()

$renewData = @(
  @(0, 20202100034, 2591.333333333333, 0, "OtherPV"),
  @(1, 20152100032, 18148.27119466832, 0, "OtherPV"),
  @(2, 20152100043, 2591.333333333333, 0, "OtherPV"),
  @(3, 20112100035, 17185.46324999998, 0, "OtherPV"),
  @(4, 20102100024, 3232.501133333333, 0, "OtherPV"),
  @(5, 20102100042, 18148.27119466832, 0, "OtherPV"),
  @(6, 20062100023, 17185.46324999998, 0, "OtherPV"),
  @(7, 20052100021, 3232.501133333333, 0, "OtherPV"),
  @(8, 20002100036, 3232.501133333333, 0, "OtherPV")
)

for ($i = 0; $i -lt $renewData.Count; $i++) {
    $r = $i + 2
    $row = $renewData[$i]
    $wsRenew.Cells.Item($r, 1).Value = $row[0]
    $wsRenew.Cells.Item($r, 2).Value = $row[1]
    $wsRenew.Cells.Item($r, 3).Value = $row[2]
    $wsRenew.Cells.Item($r, 4).Value = $row[3]
    $wsRenew.Cells.Item($r, 5).Value = $row[4]
}

# --- storages: drop the only data row, keep header ---
$wsStorage = $wb.Worksheets.Item("storages")
$wsStorage.Rows.Item(2).Delete()

# --- biogas: replace rows 2-4 with new data, drop row 5 ---
$wsBiogas = $wb.Worksheets.Item("biogas")

$wsBiogas.Range("A5:I5").Delete()

$biogasData = @(
  @(0, 20280100043, 100),
  @(1, 99990100004, 100),
  @(2, 20240100043, 100)
)

for ($i = 0; $i -lt $biogasData.Count; $i++) {
    $r = $i + 2
    $row = $biogasData[$i]
    $wsBiogas.Cells.Item($r, 1).Value = $row[0]
    $wsBiogas.Cells.Item($r, 2).Value = $row[1]
    $wsBiogas.Cells.Item($r, 3).Value = $row[2]
}
